# "Generate Report for Handoff"
#
# The localization handoff file names embed a source-file GUID and a
# content hash; a new handoff run changed both, and the handoff
# timestamps advanced. Propagate the new identifiers/timestamps across
# every sheet: the underlying cell text (and the shared strings backing
# it) as well as the cached hyperlink display text, which Excel stores
# independently of the cell value.

$wb = $excel.ActiveWorkbook

$oldGuid = "5a3d97a1-84ac-4f0e-bf2c-519e6b88ba60"
$newGuid = "a0584535-e400-4cae-9331-50bad476ca56"
$oldHash = "2f06c3c12ba4ed4a24f12403c13782a4a9d07fd4"
$newHash = "87f11f83a5d8c5ef728bc2541d53f46607f47b92"

# Update every cell value that carries the old GUID/hash (file-name cells)
# across all three sheets in one pass.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace($oldGuid, $newGuid)
    $ws.Cells.Replace($oldHash, $newHash)
}

# Hyperlink "display" text is independent of the cell's text value, so it
# needs to be refreshed separately.
foreach ($ws in $wb.Worksheets) {
    foreach ($hl in $ws.Hyperlinks) {
        $t = $hl.TextToDisplay
        $t = $t -replace [regex]::Escape($oldGuid), $newGuid
        $t = $t -replace [regex]::Escape($oldHash), $newHash
        $hl.TextToDisplay = $t
    }
}

# New handoff timestamps for the zh-cn and de-de targets.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-03-09 10:03:12"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-03-09 10:03:15"
